# aggiornamento fino a 1/09/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (357) down onto the
# newly appended rows (358:366) so column A keeps the date style (s="2").
$ws.Range("A357").Copy()
$ws.Range("A358:A366").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(44432, 0, 1, 145.7725947521866),
    @(44433, 0, 1, 145.7725947521866),
    @(44434, 2, 2, 291.5451895043732),
    @(44435, 1, 3, 437.3177842565598),
    @(44436, 4, 7, 1020.408163265306),
    @(44437, 0, 7, 1020.408163265306),
    @(44438, 0, 7, 1020.408163265306),
    @(44439, 0, 7, 1020.408163265306),
    @(44440, 0, 7, 1020.408163265306)
)

$r = 358
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
